$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text (it stores things like "505.08" or
# "57.024.41" as literal strings, not numbers) and keep the default "Normal"
# cell style so we do not introduce a new number format on these cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.024.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.395.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.407.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.321"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.825.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.928.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.410.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "309.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("E39").Value = "  +4.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.819"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.56%  "
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "133.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.18%  "
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "252.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.03%  "
